# Add a new row ("Part 028", 1 copy, PETG) to the "Tabella1" table on
# the active worksheet, then move the selection the way the author left
# it (cell E14) after scrolling the sheet down a bit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data lives in an Excel Table ("ListObject"); adding a ListRow is
# the natural way to grow it — Excel automatically extends the table
# range / AutoFilter range and the worksheet's used-range dimension.
$table = $ws.ListObjects.Item("Tabella1")
$newRow = $table.ListRows.Add()

$newRow.Range.Item(1).Value = "Part 028"
$newRow.Range.Item(2).Value = 1
# Match the centered alignment used by the existing "Number of copies"
# column cells.
$newRow.Range.Item(2).HorizontalAlignment = -4108
$newRow.Range.Item(3).Value = "PETG"

# Reflect the author's final cursor position/selection.
$ws.Range("E14").Select() | Out-Null
